$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a paragraph's content (given by a unique anchor text found
# inside it) together with the FULL text of the immediately preceding
# paragraph. Spanning the boundary of the previous paragraph avoids a stray
# leftover <w:proofErr/> marker that the underlying engine otherwise leaves
# behind when a replacement range starts exactly at a run/markup boundary.
# ---------------------------------------------------------------------------
function Replace-ParagraphWithPrevious {
    param(
        [string]$AnchorText,
        [string]$PrevParagraphXml,
        [string]$TargetParagraphXml
    )
    # NOTE: this function must be called with *positional* arguments; named
    # arguments (-AnchorText "...") are not bound correctly by this runtime.

    $full = $d.Content
    $found = $full.Find.Execute($AnchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "ANCHOR NOT FOUND:" $AnchorText
        return
    }

    $targetPara = $full.Paragraphs(1)
    $prevPara = $targetPara.Previous()

    $rangeStart = $prevPara.Range.Start
    $rangeEnd = $targetPara.Range.End

    $r = $d.Range($rangeStart, $rangeEnd)

    $xml = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'>
<w:body>
$PrevParagraphXml
$TargetParagraphXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Edit 1 (lower in the document, applied first so earlier offsets in the
# document stay valid for Edit 2): expand the abbreviated "html, css, php,
# bootstrap e JavaScript" list into explicit "HTML, CSS, PHP, Bootstrap e
# JavaScript" runs, dropping the spell-check markers around HTML/CSS/PHP and
# keeping them only around Bootstrap and JavaScript.
# ---------------------------------------------------------------------------
$prevParaXml_2 = @'
<w:p w14:paraId="78398959" w14:textId="77777777" w:rsidR="00AE31BD" w:rsidRDefault="00AE31BD" w:rsidP="003E09F2">
  <w:pPr><w:ind w:left="426"/></w:pPr>
  <w:r><w:tab/><w:t>Implementação física do projeto, feito com MySql</w:t></w:r>
</w:p>
'@

$targetParaXml_2 = @'
<w:p w14:paraId="3AB52265" w14:textId="058FC468" w:rsidR="00AE31BD" w:rsidRDefault="00AE31BD" w:rsidP="003E09F2">
  <w:pPr><w:ind w:left="426"/></w:pPr>
  <w:r><w:tab/><w:t xml:space="preserve">Início da construção do projeto usando ferramentas web como o </w:t></w:r>
  <w:r><w:t>HTML</w:t></w:r>
  <w:r><w:t xml:space="preserve">, </w:t></w:r>
  <w:r><w:t>CSS</w:t></w:r>
  <w:r><w:t xml:space="preserve">, </w:t></w:r>
  <w:r><w:t>PHP</w:t></w:r>
  <w:r><w:t xml:space="preserve">, </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Bootstrap</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve">e </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>JavaScript</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@

Replace-ParagraphWithPrevious "html, css, php, bootstrap e JavaScript" $prevParaXml_2 $targetParaXml_2

# ---------------------------------------------------------------------------
# Edit 2 (earlier in the document): merge the "(  X ) Sim" runs into a single
# run without the gramStart/gramEnd proofreading markers, and turn the tab
# before ") Não" into two literal spaces.
# ---------------------------------------------------------------------------
$prevParaXml_1 = @'
<w:p w14:paraId="0EF8BF27" w14:textId="6A8F22C0" w:rsidR="003E09F2" w:rsidRDefault="009D50CC" w:rsidP="003E09F2">
  <w:pPr><w:ind w:left="426"/></w:pPr>
  <w:r><w:t>a) -</w:t></w:r>
  <w:r w:rsidR="003E09F2"><w:t xml:space="preserve"> O cronograma das atividades está sendo executado em compatibilidade com os objetivos, metas e etapas do Projeto?</w:t></w:r>
</w:p>
'@

$targetParaXml_1 = @'
<w:p w14:paraId="2F9E10F2" w14:textId="0D448E93" w:rsidR="003E09F2" w:rsidRDefault="003E09F2" w:rsidP="003E09F2">
  <w:pPr><w:ind w:left="426"/></w:pPr>
  <w:r><w:t>(</w:t></w:r>
  <w:r w:rsidR="009D50CC"><w:t xml:space="preserve"> X ) Sim</w:t></w:r>
  <w:r><w:tab/></w:r>
  <w:r><w:tab/><w:t>(</w:t></w:r>
  <w:r w:rsidR="00117D03"><w:t xml:space="preserve">  </w:t></w:r>
  <w:r><w:t>) Não</w:t></w:r>
</w:p>
'@

Replace-ParagraphWithPrevious "(  X ) Sim" $prevParaXml_1 $targetParaXml_1

Write-Host "Done"
